$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column headers
$ws.Range("A1").Value = "building_name"
$ws.Range("B1").Value = "building_no"
$ws.Range("C1").Value = "website:map"

# Trim the stray leading/trailing spaces around the sharepoi id in the URL column (C2:C20)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $val = $cell.Value2
    if ($val -ne $null) {
        $cell.Value = ($val -replace " ", "")
    }
}
